$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add a new worksheet "test_signal" as the LAST sheet (after "Лист3").
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ts = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ts.Name = "test_signal"

# ---------------------------------------------------------------------------
# 2. Populate "test_signal".
#    NOTE: the very first *new* shared string to be introduced must be
#    "tau ns" (it has to land on shared-string index 25), so write the Q
#    column before the I:L columns that introduce f4_good/f4_bad/f5_good/f5_bad.
# ---------------------------------------------------------------------------

# -- Row 1 / Row 13: headers ------------------------------------------------
$ts.Range("Q1").Value  = "tau ns"
$ts.Range("Q13").Value = "tau ns"

$ts.Range("B1").Value  = "all"
$ts.Range("C1").Value  = "f1_good"
$ts.Range("D1").Value  = "f1_bad"
$ts.Range("E1").Value  = "f2_good"
$ts.Range("F1").Value  = "f2_bad"
$ts.Range("G1").Value  = "f3_good"
$ts.Range("H1").Value  = "f3_bad"
$ts.Range("I1").Value  = "f4_good"
$ts.Range("J1").Value  = "f4_bad"
$ts.Range("K1").Value  = "f5_good"
$ts.Range("L1").Value  = "f5_bad"
$ts.Range("O1").Value  = "chi2_per_dof_th"

$ts.Range("B13").Value = "all"
$ts.Range("C13").Value = "f1_good"
$ts.Range("D13").Value = "f1_bad"
$ts.Range("E13").Value = "f2_good"
$ts.Range("F13").Value = "f2_bad"
$ts.Range("G13").Value = "f3_good"
$ts.Range("H13").Value = "f3_bad"
$ts.Range("I13").Value = "f4_good"
$ts.Range("J13").Value = "f4_bad"
$ts.Range("K13").Value = "f5_good"
$ts.Range("L13").Value = "f5_bad"
$ts.Range("O13").Value = "chi2_per_dof_th"

# -- Row 2: summary row -------------------------------------------------
$ts.Range("A2").Value = "всего"
$ts.Range("B2").Value = 39
$ts.Range("C2").Value = 29
$ts.Range("E2").Value = 6
$ts.Range("G2").Value = 1
$ts.Range("I2").Value = 1
$ts.Range("K2").Value = 1
$ts.Range("L2").Value = 1
$ts.Range("O2").Value = 1
$ts.Range("Q2").Value = 150

# -- Row 3: "шумы" label ------------------------------------------------
$ts.Range("A3").Value = "шумы"

# -- Rows 4-9: per-signal breakdown -------------------------------------
$ts.Range("A4").Value = 1
$ts.Range("C4").Value = 29

$ts.Range("A5").Value = 2
$ts.Range("E5").Value = 6
$ts.Range("G5").Value = 1

$ts.Range("A6").Value = 3

$ts.Range("A7").Value = 4

$ts.Range("A8").Value = 5
$ts.Range("I8").Value = 1
$ts.Range("K8").Value = 1

$ts.Range("A9").Value = 6
$ts.Range("L9").Value = 1

# -- Rows 14-16: aggregate table with chi2/dof ratios --------------------
$ts.Range("A14").Value = "всего"
$ts.Range("B14").Value = 1080
$ts.Range("C14").Value = 604
$ts.Range("E14").Value = 258
$ts.Range("G14").Value = 133
$ts.Range("I14").Value = 52
$ts.Range("K14").Value = 20
$ts.Range("N14").Formula = "=SUM(C14:K14)/B14"
$ts.Range("O14").Value = 1
$ts.Range("Q14").Value = 150

$ts.Range("B15").Value = 1080
$ts.Range("C15").Value = 596
$ts.Range("E15").Value = 246
$ts.Range("G15").Value = 134
$ts.Range("I15").Value = 53
$ts.Range("K15").Value = 29
$ts.Range("N15").Formula = "=SUM(C15:K15)/B15"
$ts.Range("O15").Value = 0.1

$ts.Range("B16").Value = 1080
$ts.Range("C16").Value = 576
$ts.Range("E16").Value = 234
$ts.Range("G16").Value = 128
$ts.Range("I16").Value = 61
$ts.Range("K16").Value = 41
$ts.Range("N16").Formula = "=SUM(C16:K16)/B16"
$ts.Range("O16").Value = 0.01

# Leave the cursor on N17, matching where the original author's selection
# ended up after finishing the table.
$ts.Range("N17").Select()

# ---------------------------------------------------------------------------
# 3. Update the "295k" sheet view: the selection moved down to A45:J50
#    while scrolled so row 34 is at the top of the viewport.
# ---------------------------------------------------------------------------
$ws295 = $wb.Worksheets.Item("295k")
$ws295.Activate()
$ws295.Range("A34").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws295.Range("A45:J50").Select()

# ---------------------------------------------------------------------------
# 4. "test_signal" is the sheet the author ended up looking at, so make it
#    the active sheet/tab again (matches activeTab moving to the new sheet).
# ---------------------------------------------------------------------------
$ts.Activate()
$ts.Range("N17").Select()
